# Generate Report for Handoff
# Updates the localization-status report to reflect that the handoff
# report has just been (re)generated: status moves from "In Translation"
# to "Ready for handoff", timestamps are refreshed, and the "Status" /
# per-language status columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text from "In Translation" to "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refresh the handoff generation timestamps ---
$wsOverview.Range("G2").Value = "2016-08-20 06:46:06"
$wsDeDe.Range("H2").Value = "2016-08-20 06:46:06"
$wsZhCn.Range("H2").Value = "2016-08-20 06:45:58"

# --- Widen the Status columns to fit the new, longer text ---
# (target OOXML column width is 17.2159881591797 characters; COM's
# ColumnWidth quantizes to the screen pixel grid, so we feed the value
# that lands on the nearest representable pixel width)
$wsOverview.Range("E:E").ColumnWidth = 16.33333333333333
$wsOverview.Range("F:F").ColumnWidth = 16.33333333333333
$wsZhCn.Range("C:C").ColumnWidth = 16.33333333333333
$wsDeDe.Range("C:C").ColumnWidth = 16.33333333333333
